$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price column (D) values. These are stored as text strings that
# look numeric, so we use a leading apostrophe to force text storage and
# then restore the Normal style so no stray formatting/quote-prefix remains.
$ws.Range("D2").Value = "'243.52"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'23.69"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'5.255"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.05806"
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").Value = "'3.330"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.8082"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.8752"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'0.1386"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.07274"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'0.03064"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'0.03053"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'0.09306"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'3.852"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'0.001551"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'0.04705"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'0.0006039"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'0.006091"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'0.001267"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").Value = "'0.00008698"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").Value = "'2.144"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'0.3210"
$ws.Range("D25").Style = "Normal"
$ws.Range("D28").Value = "'0.0002343"
$ws.Range("D28").Style = "Normal"
$ws.Range("D40").Value = "'0.03789"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'0.006328"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Value = "'0.1055"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'0.002459"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'0.006956"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005486"
$ws.Range("D45").Style = "Normal"
$ws.Range("D48").Value = "'0.006926"
$ws.Range("D48").Style = "Normal"

# Update the symbol/description text in column E for rows 43 and 44
$ws.Range("E43").Value = "42CEJICEJIWorstin24h"
$ws.Range("E44").Value = "43LocalTradersLCT"
